# Auto-generated script: insert the "Knärot" section before the
# "Tretåig hackspett – ekologi samt krav på livsmiljön" Heading 1 paragraph,
# and bump the header date from 2023-09-13 to 2023-09-15.

$d = $word.ActiveDocument

# Locate the anchor paragraph: "Tretåig hackspett – ekologi samt krav på livsmiljön"
# (the first Heading 1 paragraph whose text starts with "Tretåig hackspett").
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Tretåig hackspett – ekologi") -and $p.Style.NameLocal -eq "Heading 1") {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq 0) { throw "Anchor paragraph not found" }

# --- New paragraph 0 ---
$anchor = $d.Paragraphs.Item($anchorIndex).Range
$anchor.Collapse(1)
[void]$anchor.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($anchorIndex)
$newPara.Style = "Heading 1"
$anchorIndex = $anchorIndex + 1

$pRange = $newPara.Range
$pRange.InsertAfter('Knärot – ekologi samt krav på livsmiljön')
$pRange = $newPara.Range

# --- New paragraph 1 ---
$anchor = $d.Paragraphs.Item($anchorIndex).Range
$anchor.Collapse(1)
[void]$anchor.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($anchorIndex)
$newPara.Style = "Normal"
$anchorIndex = $anchorIndex + 1

$pRange = $newPara.Range
$pRange.InsertAfter('Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).')
$pRange = $newPara.Range

# --- New paragraph 2 ---
$anchor = $d.Paragraphs.Item($anchorIndex).Range
$anchor.Collapse(1)
[void]$anchor.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($anchorIndex)
$newPara.Style = "Normal"
$anchorIndex = $anchorIndex + 1

$pRange = $newPara.Range
$pRange.InsertAfter('Samuel Johnsons doktorsavhandling ')
$pRange = $newPara.Range
$insPoint = $d.Range($pRange.End - 1, $pRange.End - 1)
$runStart = $insPoint.Start
$insPoint.InsertAfter('“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“')
$runEnd = $runStart + 82
$runRange = $d.Range($runStart, $runEnd)
$runRange.Font.Italic = 1
$pRange = $newPara.Range
$insPoint = $d.Range($pRange.End - 1, $pRange.End - 1)
$runStart = $insPoint.Start
$insPoint.InsertAfter(' (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: ')
$runEnd = $runStart + 162
$pRange = $newPara.Range
$insPoint = $d.Range($pRange.End - 1, $pRange.End - 1)
$runStart = $insPoint.Start
$insPoint.InsertAfter('“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” ')
$runEnd = $runStart + 205
$runRange = $d.Range($runStart, $runEnd)
$runRange.Font.Italic = 1
$pRange = $newPara.Range
$insPoint = $d.Range($pRange.End - 1, $pRange.End - 1)
$runStart = $insPoint.Start
$insPoint.InsertAfter('Vidare ')
$runEnd = $runStart + 7
$pRange = $newPara.Range
$insPoint = $d.Range($pRange.End - 1, $pRange.End - 1)
$runStart = $insPoint.Start
$insPoint.InsertAfter('“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”')
$runEnd = $runStart + 118
$runRange = $d.Range($runStart, $runEnd)
$runRange.Font.Italic = 1
$pRange = $newPara.Range

# --- New paragraph 3 ---
$anchor = $d.Paragraphs.Item($anchorIndex).Range
$anchor.Collapse(1)
[void]$anchor.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($anchorIndex)
$newPara.Style = "Normal"
$anchorIndex = $anchorIndex + 1

$pRange = $newPara.Range
$pRange.InsertAfter('Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: ')
$pRange = $newPara.Range
$insPoint = $d.Range($pRange.End - 1, $pRange.End - 1)
$runStart = $insPoint.Start
$insPoint.InsertAfter('“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”')
$runEnd = $runStart + 865
$runRange = $d.Range($runStart, $runEnd)
$runRange.Font.Italic = 1
$pRange = $newPara.Range

# --- New paragraph 4 ---
$anchor = $d.Paragraphs.Item($anchorIndex).Range
$anchor.Collapse(1)
[void]$anchor.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($anchorIndex)
$newPara.Style = "Normal"
$anchorIndex = $anchorIndex + 1

$pRange = $newPara.Range
$pRange.InsertAfter('En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).')
$pRange = $newPara.Range

# --- New paragraph 5 ---
$anchor = $d.Paragraphs.Item($anchorIndex).Range
$anchor.Collapse(1)
[void]$anchor.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($anchorIndex)
$newPara.Style = "Normal"
$anchorIndex = $anchorIndex + 1

$pRange = $newPara.Range
$pRange.InsertAfter('Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).')
$pRange = $newPara.Range

# --- New paragraph 6 ---
$anchor = $d.Paragraphs.Item($anchorIndex).Range
$anchor.Collapse(1)
[void]$anchor.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($anchorIndex)
$newPara.Style = "Heading 2"
$anchorIndex = $anchorIndex + 1

$pRange = $newPara.Range
$pRange.InsertAfter('Referenser - knärot')
$pRange = $newPara.Range

# --- New paragraph 7 ---
$anchor = $d.Paragraphs.Item($anchorIndex).Range
$anchor.Collapse(1)
[void]$anchor.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($anchorIndex)
$newPara.Style = "Normal"
$anchorIndex = $anchorIndex + 1

$pRange = $newPara.Range
$pRange.InsertAfter('de Graaf M & Roberts M.R., 2009. ')
$pRange = $newPara.Range
$insPoint = $d.Range($pRange.End - 1, $pRange.End - 1)
$runStart = $insPoint.Start
$insPoint.InsertAfter('Short-term response of the herbaceous layer within leave patches after harvest. ')
$runEnd = $runStart + 80
$runRange = $d.Range($runStart, $runEnd)
$runRange.Font.Italic = 1
$pRange = $newPara.Range
$insPoint = $d.Range($pRange.End - 1, $pRange.End - 1)
$runStart = $insPoint.Start
$insPoint.InsertAfter('Forest Ecology and Management 257, 1014-1025')
$runEnd = $runStart + 44
$pRange = $newPara.Range

# --- New paragraph 8 ---
$anchor = $d.Paragraphs.Item($anchorIndex).Range
$anchor.Collapse(1)
[void]$anchor.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($anchorIndex)
$newPara.Style = "Normal"
$anchorIndex = $anchorIndex + 1

$pRange = $newPara.Range
$pRange.InsertAfter('Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. ')
$pRange = $newPara.Range
$insPoint = $d.Range($pRange.End - 1, $pRange.End - 1)
$runStart = $insPoint.Start
$insPoint.InsertAfter('Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. ')
$runEnd = $runStart + 114
$runRange = $d.Range($runStart, $runEnd)
$runRange.Font.Italic = 1
$pRange = $newPara.Range
$insPoint = $d.Range($pRange.End - 1, $pRange.End - 1)
$runStart = $insPoint.Start
$insPoint.InsertAfter('Ecological Applications, 22, 2049-2064 ')
$runEnd = $runStart + 39
$pRange = $newPara.Range

# --- New paragraph 9 ---
$anchor = $d.Paragraphs.Item($anchorIndex).Range
$anchor.Collapse(1)
[void]$anchor.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($anchorIndex)
$newPara.Style = "Normal"
$anchorIndex = $anchorIndex + 1

$pRange = $newPara.Range
$pRange.InsertAfter('Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. ')
$pRange = $newPara.Range
$insPoint = $d.Range($pRange.End - 1, $pRange.End - 1)
$runStart = $insPoint.Start
$insPoint.InsertAfter('Interactive effects of drought and edge exposure on old-growth forest understory species. ')
$runEnd = $runStart + 90
$runRange = $d.Range($runStart, $runEnd)
$runRange.Font.Italic = 1
$pRange = $newPara.Range
$insPoint = $d.Range($pRange.End - 1, $pRange.End - 1)
$runStart = $insPoint.Start
$insPoint.InsertAfter('Landscape Ecology, 37, sid 1839-1853')
$runEnd = $runStart + 36
$pRange = $newPara.Range

# --- New paragraph 10 ---
$anchor = $d.Paragraphs.Item($anchorIndex).Range
$anchor.Collapse(1)
[void]$anchor.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($anchorIndex)
$newPara.Style = "Normal"
$anchorIndex = $anchorIndex + 1

$pRange = $newPara.Range
$pRange.InsertAfter('Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. ')
$pRange = $newPara.Range
$insPoint = $d.Range($pRange.End - 1, $pRange.End - 1)
$runStart = $insPoint.Start
$insPoint.InsertAfter('Biological legacies buffer local species extinction after logging. ')
$runEnd = $runStart + 67
$runRange = $d.Range($runStart, $runEnd)
$runRange.Font.Italic = 1
$pRange = $newPara.Range
$insPoint = $d.Range($pRange.End - 1, $pRange.End - 1)
$runStart = $insPoint.Start
$insPoint.InsertAfter('Journal of Applied Ecology. 51, 53-62.')
$runEnd = $runStart + 38
$pRange = $newPara.Range

# --- New paragraph 11 ---
$anchor = $d.Paragraphs.Item($anchorIndex).Range
$anchor.Collapse(1)
[void]$anchor.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($anchorIndex)
$newPara.Style = "Normal"
$anchorIndex = $anchorIndex + 1

$pRange = $newPara.Range
$pRange.InsertAfter('Skogsstyrelsen, 2022. ')
$pRange = $newPara.Range
$insPoint = $d.Range($pRange.End - 1, $pRange.End - 1)
$runStart = $insPoint.Start
$insPoint.InsertAfter('Vägledning för hänsyn till knärot. ')
$runEnd = $runStart + 35
$runRange = $d.Range($runStart, $runEnd)
$runRange.Font.Italic = 1
$pRange = $newPara.Range
$insPoint = $d.Range($pRange.End - 1, $pRange.End - 1)
$runStart = $insPoint.Start
$insPoint.InsertAfter('https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/')
$runEnd = $runStart + 128
$pRange = $newPara.Range

# --- New paragraph 12 ---
$anchor = $d.Paragraphs.Item($anchorIndex).Range
$anchor.Collapse(1)
[void]$anchor.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($anchorIndex)
$newPara.Style = "Normal"
$anchorIndex = $anchorIndex + 1

$pRange = $newPara.Range
$pRange.InsertAfter('SLU Artdatabanken, 2021. ')
$pRange = $newPara.Range
$insPoint = $d.Range($pRange.End - 1, $pRange.End - 1)
$runStart = $insPoint.Start
$insPoint.InsertAfter('Artfaktablad. Naturvård – artfakta. ')
$runEnd = $runStart + 36
$runRange = $d.Range($runStart, $runEnd)
$runRange.Font.Italic = 1
$pRange = $newPara.Range
$insPoint = $d.Range($pRange.End - 1, $pRange.End - 1)
$runStart = $insPoint.Start
$insPoint.InsertAfter('SLU Artdatabanken, Uppsala ')
$runEnd = $runStart + 27
$pRange = $newPara.Range

# --- Update header date 2023-09-13 -> 2023-09-15 ---
$d.Content.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2) | Out-Null

foreach ($sec in $d.Sections) {
    foreach ($hdrType in 1,2,3) {
        $hdr = $sec.Headers.Item($hdrType)
        if ($hdr.Exists) {
            $hdr.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2) | Out-Null
        }
    }
}

Write-Output "done"
